$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("D5").Value = 2.469545652167673
$ws.Range("D6").Value = 0.0707115266285713
$ws.Range("D7").Value = -0.3610816495766939
$ws.Range("D8").Value = 0.2166748852961936
$ws.Range("D9").Value = 2.472840372559594
$ws.Range("D10").Value = 0.2886132232229098
$ws.Range("D11").Value = 2.433730285853511
$ws.Range("D12").Value = 0.01801372628179249
$ws.Range("D13").Value = 0.3368671780379491
$ws.Range("D14").Value = 0.3687202260799354
$ws.Range("D15").Value = 0.256398869600323
$ws.Range("D16").Value = 0.215621245755459
$ws.Range("D17").Value = 0.1398250041155577
$ws.Range("D18").Value = -0.01753936205462238
$ws.Range("D19").Value = 0.003338890064221986
$ws.Range("D20").Value = 0.4428044280442804
$ws.Range("D21").Value = 0.0005665001879956587
$ws.Range("D22").Value = 0.4818836742598442
$ws.Range("D23").Value = 0.2404175451765602
